# edit.ps1 -- Cam7calib.xlsx calibration-data refresh
# Rewrites the A:D calibration table on Sheet1 (rows 1-103) with updated
# measurement values, blanks out the now-unused rows 104-108, trims the
# trailing empty rows 268-276 off the bottom of the sheet, and leaves the
# selection on D83 (the cell last edited by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the calibration values in A1:D103 -------------------------
$data = New-Object 'object[,]' 103,4

$data[0,0] = 929
$data[0,1] = 1030.656
$data[0,2] = 2500.5
$data[0,3] = 0
$data[1,0] = 951.607
$data[1,1] = 1804.086
$data[1,2] = 2500.5
$data[1,3] = -16.08
$data[2,0] = 945.943
$data[2,1] = 1700.328
$data[2,2] = 2500.5
$data[2,3] = -14
$data[3,0] = 941.218
$data[3,1] = 1600.713
$data[3,2] = 2500.5
$data[3,3] = -12
$data[4,0] = 937.237
$data[4,1] = 1503.705
$data[4,2] = 2500.5
$data[4,3] = -9.98
$data[5,0] = 934.136
$data[5,1] = 1408.869
$data[5,2] = 2500.5
$data[5,3] = -7.98
$data[6,0] = 932.228
$data[6,1] = 1332.617
$data[6,2] = 2500.5
$data[6,3] = -5.98
$data[7,0] = 930.038
$data[7,1] = 1219.703
$data[7,2] = 2500.5
$data[7,3] = -4
$data[8,0] = 929.021
$data[8,1] = 1031.083
$data[8,2] = 2500.5
$data[8,3] = -0.03
$data[9,0] = 929.509
$data[9,1] = 935.722
$data[9,2] = 2500.5
$data[9,3] = 2
$data[10,0] = 930.791
$data[10,1] = 840.975
$data[10,2] = 2500.5
$data[10,3] = 3.97
$data[11,0] = 932.549
$data[11,1] = 744.609
$data[11,2] = 2500.5
$data[11,3] = 6
$data[12,0] = 935.389
$data[12,1] = 648.696
$data[12,2] = 2500.5
$data[12,3] = 7.95
$data[13,0] = 938.583
$data[13,1] = 549.688
$data[13,2] = 2500.5
$data[13,3] = 9.97
$data[14,0] = 942.698
$data[14,1] = 451.139
$data[14,2] = 2500.5
$data[14,3] = 11.97
$data[15,0] = 947.57
$data[15,1] = 350.789
$data[15,2] = 2500.5
$data[15,3] = 13.97
$data[16,0] = 953.445
$data[16,1] = 248.605
$data[16,2] = 2500.5
$data[16,3] = 15.95
$data[17,0] = 910.097
$data[17,1] = 1030.947
$data[17,2] = 2600.2
$data[17,3] = 0
$data[18,0] = 932.155
$data[18,1] = 1802.521
$data[18,2] = 2600.2
$data[18,3] = -16.08
$data[19,0] = 926.548
$data[19,1] = 1699.31
$data[19,2] = 2600.2
$data[19,3] = -14
$data[20,0] = 922.096
$data[20,1] = 1599.967
$data[20,2] = 2600.2
$data[20,3] = -11.98
$data[21,0] = 918.329
$data[21,1] = 1504.07
$data[21,2] = 2600.2
$data[21,3] = -9.98
$data[22,0] = 915.242
$data[22,1] = 1409.221
$data[22,2] = 2600.2
$data[22,3] = -8
$data[23,0] = 913.002
$data[23,1] = 1313.242
$data[23,2] = 2600.2
$data[23,3] = -5.98
$data[24,0] = 911.38
$data[24,1] = 1219.485
$data[24,2] = 2600.2
$data[24,3] = -4
$data[25,0] = 910.555
$data[25,1] = 1124.303
$data[25,2] = 2600.2
$data[25,3] = -1.98
$data[26,0] = 910.107
$data[26,1] = 1031.193
$data[26,2] = 2600.2
$data[26,3] = 0
$data[27,0] = 910.97
$data[27,1] = 936.838
$data[27,2] = 2600.2
$data[27,3] = 1.97
$data[28,0] = 911.585
$data[28,1] = 859.902
$data[28,2] = 2600.2
$data[28,3] = 4
$data[29,0] = 913.644
$data[29,1] = 745.661
$data[29,2] = 2600.2
$data[29,3] = 5.97
$data[30,0] = 916.421
$data[30,1] = 648.904
$data[30,2] = 2600.2
$data[30,3] = 7.97
$data[31,0] = 869.143
$data[31,1] = 550.962
$data[31,2] = 2600.2
$data[31,3] = 9.95
$data[32,0] = 923.588
$data[32,1] = 450.94
$data[32,2] = 2600.2
$data[32,3] = 11.97
$data[33,0] = 928.416
$data[33,1] = 351.058
$data[33,2] = 2600.2
$data[33,3] = 13.97
$data[34,0] = 933.821
$data[34,1] = 248.717
$data[34,2] = 2600.2
$data[34,3] = 15.95
$data[35,0] = 893.001
$data[35,1] = 1030.657
$data[35,2] = 2699.9
$data[35,3] = 0
$data[36,0] = 914.066
$data[36,1] = 1801.834
$data[36,2] = 2699.9
$data[36,3] = -16.08
$data[37,0] = 908.673
$data[37,1] = 1698.238
$data[37,2] = 2699.9
$data[37,3] = -14
$data[38,0] = 900.625
$data[38,1] = 1502.662
$data[38,2] = 2699.9
$data[38,3] = -9.98
$data[39,0] = 897.952
$data[39,1] = 1407.484
$data[39,2] = 2699.9
$data[39,3] = -7.98
$data[40,0] = 895.849
$data[40,1] = 1312.899
$data[40,2] = 2699.9
$data[40,3] = -5.98
$data[41,0] = 894.003
$data[41,1] = 1219.84
$data[41,2] = 2699.9
$data[41,3] = -4
$data[42,0] = 893.104
$data[42,1] = 1125.193
$data[42,2] = 2699.9
$data[42,3] = -1.98
$data[43,0] = 892.981
$data[43,1] = 1031.902
$data[43,2] = 2699.9
$data[43,3] = -0.03
$data[44,0] = 893.115
$data[44,1] = 936.815
$data[44,2] = 2699.9
$data[44,3] = 2
$data[45,0] = 894.338
$data[45,1] = 841.698
$data[45,2] = 2699.9
$data[45,3] = 4
$data[46,0] = 896.195
$data[46,1] = 746.33
$data[46,2] = 2699.9
$data[46,3] = 5.97
$data[47,0] = 905.907
$data[47,1] = 452.37
$data[47,2] = 2699.9
$data[47,3] = 11.97
$data[48,0] = 910.357
$data[48,1] = 352.797
$data[48,2] = 2699.9
$data[48,3] = 13.97
$data[49,0] = 915.917
$data[49,1] = 249.981
$data[49,2] = 2699.9
$data[49,3] = 15.95
$data[50,0] = 876.942
$data[50,1] = 1030.75
$data[50,2] = 2799.6
$data[50,3] = 0
$data[51,0] = 897.464
$data[51,1] = 1800.667
$data[51,2] = 2799.6
$data[51,3] = -16.08
$data[52,0] = 892.375
$data[52,1] = 1697.575
$data[52,2] = 2799.6
$data[52,3] = -14
$data[53,0] = 887.954
$data[53,1] = 1599.728
$data[53,2] = 2799.6
$data[53,3] = -11.98
$data[54,0] = 884.5
$data[54,1] = 1502.911
$data[54,2] = 2799.6
$data[54,3] = -9.98
$data[55,0] = 881.829
$data[55,1] = 1408.301
$data[55,2] = 2799.6
$data[55,3] = -8
$data[56,0] = 879.364
$data[56,1] = 1313.012
$data[56,2] = 2799.6
$data[56,3] = -5.98
$data[57,0] = 877.99
$data[57,1] = 1219.47
$data[57,2] = 2799.6
$data[57,3] = -4
$data[58,0] = 877.049
$data[58,1] = 1124.576
$data[58,2] = 2799.6
$data[58,3] = -1.98
$data[59,0] = 876.915
$data[59,1] = 1031.547
$data[59,2] = 2799.6
$data[59,3] = 0
$data[60,0] = 877.108
$data[60,1] = 937.547
$data[60,2] = 2799.6
$data[60,3] = 1.97
$data[61,0] = 878.067
$data[61,1] = 841.763
$data[61,2] = 2799.6
$data[61,3] = 4
$data[62,0] = 879.816
$data[62,1] = 766.288
$data[62,2] = 2799.6
$data[62,3] = 5.97
$data[63,0] = 882.536
$data[63,1] = 650.092
$data[63,2] = 2799.6
$data[63,3] = 7.97
$data[64,0] = 885.76
$data[64,1] = 552.281
$data[64,2] = 2799.6
$data[64,3] = 9.95
$data[65,0] = 889.352
$data[65,1] = 452.289
$data[65,2] = 2799.6
$data[65,3] = 11.97
$data[66,0] = 893.858
$data[66,1] = 353.352
$data[66,2] = 2799.6
$data[66,3] = 13.97
$data[67,0] = 899.027
$data[67,1] = 250.74
$data[67,2] = 2799.6
$data[67,3] = 15.95
$data[68,0] = 861.935
$data[68,1] = 1030.848
$data[68,2] = 2899.3
$data[68,3] = 0
$data[69,0] = 881.623
$data[69,1] = 1800.321
$data[69,2] = 2899.3
$data[69,3] = -16.07
$data[70,0] = 876.723
$data[70,1] = 1696.662
$data[70,2] = 2899.3
$data[70,3] = -14
$data[71,0] = 872.512
$data[71,1] = 1598.646
$data[71,2] = 2899.3
$data[71,3] = -12
$data[72,0] = 869.301
$data[72,1] = 1501.234
$data[72,2] = 2899.3
$data[72,3] = -9.97
$data[73,0] = 866.5
$data[73,1] = 1406.266
$data[73,2] = 2899.3
$data[73,3] = -7.97
$data[74,0] = 864.374
$data[74,1] = 1312.108
$data[74,2] = 2899.3
$data[74,3] = -6
$data[75,0] = 862.071
$data[75,1] = 1124.541
$data[75,2] = 2899.3
$data[75,3] = -2
$data[76,0] = 861.893
$data[76,1] = 1030.553
$data[76,2] = 2899.3
$data[76,3] = 0
$data[77,0] = 862.222
$data[77,1] = 936.361
$data[77,2] = 2899.3
$data[77,3] = 2
$data[78,0] = 863.169
$data[78,1] = 842.175
$data[78,2] = 2899.3
$data[78,3] = 3.98
$data[79,0] = 865.001
$data[79,1] = 746.134
$data[79,2] = 2899.3
$data[79,3] = 5.98
$data[80,0] = 867.343
$data[80,1] = 649.806
$data[80,2] = 2899.3
$data[80,3] = 7.95
$data[81,0] = 870.568
$data[81,1] = 551.036
$data[81,2] = 2899.3
$data[81,3] = 9.98
$data[82,0] = 874.041
$data[82,1] = 452.978
$data[82,2] = 2899.3
$data[82,3] = 11.98
$data[83,0] = 878.611
$data[83,1] = 353.324
$data[83,2] = 2899.3
$data[83,3] = 13.98
$data[84,0] = 883.45
$data[84,1] = 251.269
$data[84,2] = 2899.3
$data[84,3] = 15.93
$data[85,0] = 847.961
$data[85,1] = 1030.728
$data[85,2] = 2998.9
$data[85,3] = 0
$data[86,0] = 867.469
$data[86,1] = 1799.584
$data[86,2] = 2998.9
$data[86,3] = -16.08
$data[87,0] = 862.509
$data[87,1] = 1696.041
$data[87,2] = 2998.9
$data[87,3] = -14
$data[88,0] = 858.493
$data[88,1] = 1598.392
$data[88,2] = 2998.9
$data[88,3] = -11.98
$data[89,0] = 855.317
$data[89,1] = 1502.56
$data[89,2] = 2998.9
$data[89,3] = -9.98
$data[90,0] = 852.5
$data[90,1] = 1407.996
$data[90,2] = 2998.9
$data[90,3] = -7.98
$data[91,0] = 850.48
$data[91,1] = 1313.583
$data[91,2] = 2998.9
$data[91,3] = -6
$data[92,0] = 848.971
$data[92,1] = 1219.056
$data[92,2] = 2998.9
$data[92,3] = -3.98
$data[93,0] = 848.194
$data[93,1] = 1125.818
$data[93,2] = 2998.9
$data[93,3] = -2
$data[94,0] = 847.936
$data[94,1] = 1032.038
$data[94,2] = 2998.9
$data[94,3] = 0
$data[95,0] = 848.301
$data[95,1] = 937.878
$data[95,2] = 2998.9
$data[95,3] = 2
$data[96,0] = 849.171
$data[96,1] = 843.501
$data[96,2] = 2998.9
$data[96,3] = 3.97
$data[97,0] = 851
$data[97,1] = 747.565
$data[97,2] = 2998.9
$data[97,3] = 5.97
$data[98,0] = 853.129
$data[98,1] = 650.731
$data[98,2] = 2998.9
$data[98,3] = 7.97
$data[99,0] = 855.788
$data[99,1] = 572.944
$data[99,2] = 2998.9
$data[99,3] = 9.97
$data[100,0] = 859.787
$data[100,1] = 455.602
$data[100,2] = 2998.9
$data[100,3] = 11.95
$data[101,0] = 864.004
$data[101,1] = 354.025
$data[101,2] = 2998.9
$data[101,3] = 14
$data[102,0] = 867.949
$data[102,1] = 273.424
$data[102,2] = 2998.9
$data[102,3] = 15.95

$ws.Range("A1:D103").Value = $data

# --- 2. Rows 104-108 no longer hold data: clear them to blank cells -------
$ws.Range("A104:D108").ClearContents()

# --- 3. Drop the trailing empty rows 268-276 (sheet dimension D276->D267) -
$ws.Range("A268:D276").EntireRow.Delete() | Out-Null

# --- 4. Restore the author's final selection/scroll position --------------
$ws.Range("D83").Select() | Out-Null

Write-Output "edit applied"
